$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.025.39'
$ws.Range("E2").Value = '  -2.84%  '

$ws.Range("D3").Value = '1.707.54'
$ws.Range("E3").Value = '  -3.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.15%  '

$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4723'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3419'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.05'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07253'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.031'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.78%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.820'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.03%  '

$ws.Range("D15").Value = '1.705.99'
$ws.Range("E15").Value = '  -3.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.807'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.68'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001032'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06356'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.85%  '

$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.586'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.86%  '

$ws.Range("D23").Value = '27.056.64'
$ws.Range("E23").Value = '  -2.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.100'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.30%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.35%  '

$ws.Range("D28").Value = '1.900.21'
$ws.Range("E28").Value = '  -3.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.056'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.005'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09135'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.580'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.271'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02179'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05804'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.86%  '

$ws.Range("E37").Value = '  -6.85%  '

$ws.Range("E38").Value = '  -0.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1978'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.714'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.390'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5860'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.097'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.432'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.554'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5610'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '117.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.01%  '

$ws.Range("E49").Value = '  -6.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06622'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.99%  '

$ws.Range("E51").Value = '  -4.57%  '
